$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new data row (row 16): 2025/11/25, 逃离鸭科夫, 1256 — continuing
# the existing daily mod-count log.
#
# Column A holds dates stored as literal text (matching every prior row,
# which is t="inlineStr"/shared-string text, not a real Excel date serial).
# A plain .Value assignment of "2025/11/25" would be auto-parsed into a date
# serial by Excel's smart entry, so the cell is pre-formatted as Text ("@")
# before the value is typed in to keep it literal.
$ws.Range("A16").NumberFormat = "@"
$ws.Range("A16").Value = "2025/11/25"
$ws.Range("B16").Value = "逃离鸭科夫"
$ws.Range("C16").Value = 1256

# Match the formatting (centered alignment, General number format) used by
# the rest of the data rows by copying row 15's formats onto row 16. Doing
# this *after* typing the value also cleans up the temporary "@" text
# format applied to A16 above, so the final style matches the other rows
# exactly.
$ws.Range("A15:C15").Copy()
$ws.Range("A16:C16").PasteSpecial(-4122)
